# Update column F (dSF) values to reflect the repulled / recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 2
    4  = 1
    5  = 1
    6  = -5
    7  = -2
    8  = -5
    9  = 2
    10 = 2
    11 = 1
    13 = -3
    14 = -2
    15 = -1
    16 = 3
    17 = 3
    18 = 3
    19 = -3
    20 = -3
    21 = 2
    22 = 2
    24 = -5
    25 = 3
    26 = 2
    28 = 3
    29 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
